# The dataset is a weekly price series. A new weekly observation was added
# as a new row right after the current row 7 (i.e. it becomes the new row 8),
# pushing every subsequent row down by one (old row 136 becomes row 137).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 8 - this shifts rows 8..136 down to 9..137.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly observation.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44496
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112017
$ws.Range("G8").Value = "Apio"
$ws.Range("H8").Value = "Americana (o)"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 8500
$ws.Range("M8").Value = 8250
$ws.Range("N8").Value = "`$/docena de matas"
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 1375
$ws.Range("Q8").Value = 6
$ws.Range("R8").Value = "Hortaliza"
